$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.068.93'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.892.32'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''307.10'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '''0.5177'
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("D8").Value = '''0.3755'
$ws.Range("D9").Value = '''0.07214'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").Value = '''0.8973'
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '''0.07656'
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '1.876.37'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '''94.20'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '''5.227'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '''0.000008513'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '''14.42'
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").Value = '''0.9998'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = '27.123.56'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '''5.060'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").Value = '2.121.06'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("E23").Value = '  +1.71%  '
$ws.Range("D24").Value = '''6.409'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = '''2.282'
$ws.Range("E25").Value = '  +9.78%  '
$ws.Range("D26").Value = '''146.40'
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").Value = '''1.737'
$ws.Range("E27").Value = '  -3.10%  '
$ws.Range("D28").Value = '''18.05'
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("D29").Value = '''114.39'
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").Value = '''4.981'
$ws.Range("E30").Value = '  +6.25%  '
$ws.Range("D31").Value = '''4.781'
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").Value = '''0.09193'
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").Value = '''0.05045'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("D34").Value = '''1.236'
$ws.Range("E34").Value = '  +6.61%  '
$ws.Range("D35").Value = '''0.7742'
$ws.Range("E35").Value = '  +2.93%  '
$ws.Range("D36").Value = '''2.977'
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("D37").Value = '''3.281'
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("D38").Value = '''2.592'
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").Value = '''8.984'
$ws.Range("E42").Value = '  +5.33%  '
$ws.Range("D43").Value = '''6.627'
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = '''119.22'
$ws.Range("E44").Value = '  +3.17%  '
$ws.Range("D45").Value = '''0.1516'
$ws.Range("E45").Value = '  +2.65%  '
$ws.Range("D46").Value = '''0.4817'
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '''0.9999'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''10.13'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").Value = '''1.594'
$ws.Range("E49").Value = '  +1.81%  '
$ws.Range("D50").Value = '''37.40'
$ws.Range("E50").Value = '  +1.41%  '
